$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.273.86"
$ws.Range("E2").Value = "  +2.87%  "

# Row 3
$ws.Range("D3").Value = "2.305.07"
$ws.Range("E3").Value = "  +2.10%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.67"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.66"
$ws.Range("E6").Value = "  +6.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.532"
$ws.Range("E7").Value = "  +1.75%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("E9").Value = "  +8.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  +3.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.75"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13
$ws.Range("E13").Value = "  -1.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.02"
$ws.Range("E14").Value = "  +2.88%  "

# Row 15
$ws.Range("D15").Value = "2.663.80"
$ws.Range("E15").Value = "  +2.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.99"
$ws.Range("E16").Value = "  +2.62%  "

# Row 17
$ws.Range("D17").Value = "2.276.81"
$ws.Range("E17").Value = "  +1.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.807"
$ws.Range("E18").Value = "  +2.73%  "

# Row 19
$ws.Range("D19").Value = "43.174.40"
$ws.Range("E19").Value = "  +3.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.31"
$ws.Range("E20").Value = "  +0.80%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0929"
$ws.Range("E21").Value = "  +3.39%  "

# Row 22
$ws.Range("E22").Value = "  +3.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.03"
$ws.Range("E23").Value = "  +0.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.27"
$ws.Range("E24").Value = "  +2.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  +1.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.61"
$ws.Range("E26").Value = "  +1.43%  "

# Row 27
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.78"
$ws.Range("E28").Value = "  +5.79%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.62"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  +3.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.62"
$ws.Range("E31").Value = "  +1.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.48"
$ws.Range("E32").Value = "  +2.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.27"
$ws.Range("E33").Value = "  +0.99%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.54"
$ws.Range("E35").Value = "  +6.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.86"
$ws.Range("E36").Value = "  +1.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0741"
$ws.Range("E37").Value = "  +1.20%  "

# Row 38
$ws.Range("E38").Value = "  -2.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  +4.19%  "

# Row 40
$ws.Range("E40").Value = "  +1.39%  "

# Row 41
$ws.Range("E41").Value = "  +1.58%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.33"
$ws.Range("E42").Value = "  +5.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0293"
$ws.Range("E44").Value = "  +5.14%  "

# Row 45
$ws.Range("D45").Value = "1.976.32"
$ws.Range("E45").Value = "  +1.31%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.28"
$ws.Range("E46").Value = "  +1.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.99"
$ws.Range("E47").Value = "  +3.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.01"
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("E49").Value = "  +3.38%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.36"
$ws.Range("E50").Value = "  +3.84%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.59"
$ws.Range("E51").Value = "  +9.65%  "
